$d = $word.ActiveDocument

# --- 1. Give the class writeup its real name wherever the "Some Class
#        Name" / "Class Name" placeholders appear (heading + CRC card).
#        Do the longer phrase first so it doesn't get caught by the
#        shorter "Class Name" search below. ---
$d.Content.Find.Execute("Some Class Name", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "ObserverView", 2) | Out-Null
$d.Content.Find.Execute("Class Name", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "ObserverView", 2) | Out-Null

# --- 2. The stray "_GoBack" bookmark (Word's "last edit position"
#        marker) used to sit in the middle of the description paragraph,
#        splitting "Some paragraph about this class" into two runs
#        ("So" / "me paragraph about this class"). Remove it and rejoin
#        the text into a single clean run. (Re-assigning the identical
#        text is a no-op to the engine, so nudge through a placeholder
#        first to force the runs to really be rewritten.) ---
$goBack = $d.Bookmarks("_GoBack")
$splitPos = $goBack.Range.Start
$goBack.Delete()

$joinRange = $d.Range($splitPos - 2, $splitPos + 29)
$joinRange.Text = "{{TEMP}}"
$joinRange = $d.Range($splitPos - 2, $splitPos + 6)
$joinRange.Text = "Some paragraph about this class"

# --- 3. Re-create "_GoBack" where the last edit actually happened: right
#        after the class name we just typed into the CRC card's cell. ---
$cell = $d.Tables(1).Cell(1, 1)
$cellRange = $cell.Range
$insertAt = $cellRange.End - 1

# A zero-length range placed directly at the end-of-cell boundary isn't a
# safe anchor, so nudge past it with a throwaway character, bookmark
# right before that character, then remove the character again - the
# bookmark collapses back to exactly the end of "ObserverView".
$tempRange = $d.Range($insertAt, $insertAt)
$tempRange.InsertAfter("X")

$bmRange = $d.Range($insertAt, $insertAt)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$tempCharRange = $d.Range($insertAt, $insertAt + 1)
$tempCharRange.Delete()
